$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.223.96"
$ws.Range("E2").Value = "  +1.29%  "
$ws.Range("D3").Value = "3.865.62"
$ws.Range("E3").Value = "  +1.59%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "692.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.46%  "
$ws.Range("D7").Value = "3.863.31"
$ws.Range("E7").Value = "  +1.60%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("E10").Value = "  +1.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.37"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.95%  "
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("E13").Value = "  +6.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.65"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.72%  "
$ws.Range("D15").Value = "4.517.72"
$ws.Range("E15").Value = "  +1.64%  "
$ws.Range("D16").Value = "3.868.01"
$ws.Range("E16").Value = "  +1.73%  "
$ws.Range("D17").Value = "71.280.37"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("E19").Value = "  +1.16%  "
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "493.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.58%  "
$ws.Range("E23").Value = "  +1.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000147"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.85%  "
$ws.Range("E26").Value = "  +1.64%  "
$ws.Range("E27").Value = "  +3.23%  "
$ws.Range("E28").Value = "  +1.99%  "
$ws.Range("D29").Value = "4.020.84"
$ws.Range("E29").Value = "  +1.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("E31").Value = "  +9.47%  "
$ws.Range("E32").Value = "  +3.68%  "
$ws.Range("E33").Value = "  +0.24%  "
$ws.Range("E34").Value = "  +0.71%  "
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("E36").Value = "  +2.51%  "
$ws.Range("D37").Value = "3.818.15"
$ws.Range("E37").Value = "  +1.53%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("E39").Value = "  +2.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +13.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.62%  "
$ws.Range("E42").Value = "  +2.05%  "
$ws.Range("E43").Value = "  +6.14%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("E46").Value = "  +3.15%  "
$ws.Range("E47").Value = "  +5.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.303"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.27%  "
$ws.Range("E51").Value = "  +2.28%  "
